$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# Clear the existing used range on the #system sheet before rewriting it
$ws.Range("A1:AA120").ClearContents()

$ws.Cells.Item(1,1).Value = "target"
$ws.Cells.Item(1,2).Value = "aws.s3"
$ws.Cells.Item(1,3).Value = "aws.ses"
$ws.Cells.Item(1,4).Value = "base"
$ws.Cells.Item(1,5).Value = "csv"
$ws.Cells.Item(1,6).Value = "desktop"
$ws.Cells.Item(1,7).Value = "excel"
$ws.Cells.Item(1,8).Value = "external"
$ws.Cells.Item(1,9).Value = "image"
$ws.Cells.Item(1,10).Value = "io"
$ws.Cells.Item(1,11).Value = "jms"
$ws.Cells.Item(1,12).Value = "json"
$ws.Cells.Item(1,13).Value = "macro"
$ws.Cells.Item(1,14).Value = "mail"
$ws.Cells.Item(1,15).Value = "number"
$ws.Cells.Item(1,16).Value = "pdf"
$ws.Cells.Item(1,17).Value = "rdbms"
$ws.Cells.Item(1,18).Value = "redis"
$ws.Cells.Item(1,19).Value = "sms"
$ws.Cells.Item(1,20).Value = "sound"
$ws.Cells.Item(1,21).Value = "ssh"
$ws.Cells.Item(1,22).Value = "step"
$ws.Cells.Item(1,23).Value = "web"
$ws.Cells.Item(1,24).Value = "webalert"
$ws.Cells.Item(1,25).Value = "webcookie"
$ws.Cells.Item(1,26).Value = "ws"
$ws.Cells.Item(1,27).Value = "ws.async"
$ws.Cells.Item(1,28).Value = "xml"
$ws.Cells.Item(2,1).Value = "aws.s3"
$ws.Cells.Item(2,2).Value = "assertNotPresent(profile,remotePath)"
$ws.Cells.Item(2,3).Value = "sendHtmlMail(profile,to,subject,body)"
$ws.Cells.Item(2,4).Value = "appendText(var,appendWith)"
$ws.Cells.Item(2,5).Value = "compare(expected,actual,failFast)"
$ws.Cells.Item(2,6).Value = "assertAttribute(locator,attribute,expected)"
$ws.Cells.Item(2,7).Value = "assertPassword(file)"
$ws.Cells.Item(2,8).Value = "runJUnit(className)"
$ws.Cells.Item(2,9).Value = "colorbit(source,bit,saveTo)"
$ws.Cells.Item(2,10).Value = "assertEqual(expected,actual)"
$ws.Cells.Item(2,11).Value = "receive(var,config,waitMs)"
$ws.Cells.Item(2,12).Value = "addOrReplace(json,jsonpath,input,var)"
$ws.Cells.Item(2,13).Value = "description()"
$ws.Cells.Item(2,14).Value = "send(profile,to,subject,body)"
$ws.Cells.Item(2,15).Value = "assertBetween(num,minNum,maxNum)"
$ws.Cells.Item(2,16).Value = "assertContentEqual(actualPdf,expectedPdf)"
$ws.Cells.Item(2,17).Value = "resultToCSV(var,csvFile,delim,showHeader)"
$ws.Cells.Item(2,18).Value = "append(profile,key,value)"
$ws.Cells.Item(2,19).Value = "sendText(phones,text)"
$ws.Cells.Item(2,20).Value = "laser(repeats)"
$ws.Cells.Item(2,21).Value = "scpCopyFrom(var,profile,remote,local)"
$ws.Cells.Item(2,22).Value = "observe(prompt)"
$ws.Cells.Item(2,23).Value = "assertAndClick(locator,label)"
$ws.Cells.Item(2,24).Value = "accept()"
$ws.Cells.Item(2,25).Value = "assertNotPresent(name)"
$ws.Cells.Item(2,26).Value = "assertReturnCode(var,returnCode)"
$ws.Cells.Item(2,27).Value = "delete(url,body,output)"
$ws.Cells.Item(2,28).Value = "assertCorrectness(xml,schema)"
$ws.Cells.Item(3,1).Value = "aws.ses"
$ws.Cells.Item(3,2).Value = "assertPresent(profile,remotePath)"
$ws.Cells.Item(3,3).Value = "sendTextMail(profile,to,subject,body)"
$ws.Cells.Item(3,4).Value = "assertArrayContain(array,expected)"
$ws.Cells.Item(3,5).Value = "compareExtended(var,profile,expected,actual)"
$ws.Cells.Item(3,6).Value = "assertChecked(name)"
$ws.Cells.Item(3,7).Value = "clear(file,worksheet,range)"
$ws.Cells.Item(3,8).Value = "runProgram(programPathAndParams)"
$ws.Cells.Item(3,9).Value = "compare(baseline,actual)"
$ws.Cells.Item(3,10).Value = "assertNotEqual(expected,actual)"
$ws.Cells.Item(3,11).Value = "sendMap(config,id,payload)"
$ws.Cells.Item(3,12).Value = "assertCorrectness(json,schema)"
$ws.Cells.Item(3,13).Value = "expects(var,default)"
$ws.Cells.Item(3,15).Value = "assertEqual(num1,num2)"
$ws.Cells.Item(3,16).Value = "assertFormElementPresent(var,name)"
$ws.Cells.Item(3,17).Value = "runFile(var,db,file)"
$ws.Cells.Item(3,18).Value = "assertKeyExists(profile,key)"
$ws.Cells.Item(3,20).Value = "play(audio)"
$ws.Cells.Item(3,21).Value = "scpCopyTo(var,profile,local,remote)"
$ws.Cells.Item(3,22).Value = "perform(instructions)"
$ws.Cells.Item(3,23).Value = "assertAttribute(locator,attrName,value)"
$ws.Cells.Item(3,24).Value = "assertPresent()"
$ws.Cells.Item(3,25).Value = "assertPresent(name)"
$ws.Cells.Item(3,26).Value = "delete(url,body,var)"
$ws.Cells.Item(3,27).Value = "download(url,queryString,saveTo)"
$ws.Cells.Item(3,28).Value = "assertElementCount(xml,xpath,count)"
$ws.Cells.Item(4,1).Value = "base"
$ws.Cells.Item(4,2).Value = "copyFrom(var,profile,remote,local)"
$ws.Cells.Item(4,4).Value = "assertArrayEqual(array1,array2,exactOrder)"
$ws.Cells.Item(4,5).Value = "convertExcel(excel,worksheet,csvFile)"
$ws.Cells.Item(4,6).Value = "assertDisabled(name)"
$ws.Cells.Item(4,7).Value = "clearPassword(file,password)"
$ws.Cells.Item(4,8).Value = "runProgramNoWait(programPathAndParams)"
$ws.Cells.Item(4,9).Value = "convert(source,format,saveTo)"
$ws.Cells.Item(4,10).Value = "assertReadableFile(file,minByte)"
$ws.Cells.Item(4,11).Value = "sendText(config,id,payload)"
$ws.Cells.Item(4,12).Value = "assertElementCount(json,jsonpath,count)"
$ws.Cells.Item(4,13).Value = "produces(var,value)"
$ws.Cells.Item(4,15).Value = "assertGreater(num1,num2)"
$ws.Cells.Item(4,16).Value = "assertFormValue(var,name,expected)"
$ws.Cells.Item(4,17).Value = "runSQL(var,db,sql)"
$ws.Cells.Item(4,18).Value = "delete(profile,key)"
$ws.Cells.Item(4,20).Value = "speak(text)"
$ws.Cells.Item(4,21).Value = "sftpCopyFrom(var,profile,remote,local)"
$ws.Cells.Item(4,22).Value = "validate(prompt,responses,passResponses)"
$ws.Cells.Item(4,23).Value = "assertAttributeContains(locator,attrName,contains)"
$ws.Cells.Item(4,24).Value = "assertText(text,matchBy)"
$ws.Cells.Item(4,25).Value = "assertValue(name,value)"
$ws.Cells.Item(4,26).Value = "download(url,queryString,saveTo)"
$ws.Cells.Item(4,27).Value = "get(url,queryString,output)"
$ws.Cells.Item(4,28).Value = "assertElementNotPresent(xml,xpath)"
$ws.Cells.Item(5,1).Value = "csv"
$ws.Cells.Item(5,2).Value = "copyTo(var,profile,local,remote)"
$ws.Cells.Item(5,4).Value = "assertArrayNotContain(array,unexpected)"
$ws.Cells.Item(5,5).Value = "fromExcel(excel,worksheet,csvFile)"
$ws.Cells.Item(5,6).Value = "assertElementPresent(name)"
$ws.Cells.Item(5,7).Value = "columnarCsv(file,worksheet,ranges,output)"
$ws.Cells.Item(5,9).Value = "crop(image,dimension,saveTo)"
$ws.Cells.Item(5,10).Value = "base64(var,file)"
$ws.Cells.Item(5,12).Value = "assertElementNotPresent(json,jsonpath)"
$ws.Cells.Item(5,15).Value = "assertGreaterOrEqual(num1,num2)"
$ws.Cells.Item(5,16).Value = "assertFormValues(var,name,expectedValues,exactOrder)"
$ws.Cells.Item(5,17).Value = "runSQLs(var,db,sqls)"
$ws.Cells.Item(5,18).Value = "flushAll(profile)"
$ws.Cells.Item(5,20).Value = "speakNoWait(text)"
$ws.Cells.Item(5,21).Value = "sftpCopyTo(var,profile,local,remote)"
$ws.Cells.Item(5,23).Value = "assertAttributeNotContains(locator,attrName,contains)"
$ws.Cells.Item(5,24).Value = "dismiss()"
$ws.Cells.Item(5,25).Value = "delete(name)"
$ws.Cells.Item(5,26).Value = "get(url,queryString,var)"
$ws.Cells.Item(5,27).Value = "head(url,output)"
$ws.Cells.Item(5,28).Value = "assertElementPresent(xml,xpath)"
$ws.Cells.Item(6,1).Value = "desktop"
$ws.Cells.Item(6,2).Value = "delete(var,profile,remotePath)"
$ws.Cells.Item(6,4).Value = "assertContains(text,substring)"
$ws.Cells.Item(6,6).Value = "assertEnabled(name)"
$ws.Cells.Item(6,7).Value = "csv(file,worksheet,range,output)"
$ws.Cells.Item(6,9).Value = "resize(image,width,height,saveTo)"
$ws.Cells.Item(6,10).Value = "compare(expected,actual,failFast)"
$ws.Cells.Item(6,12).Value = "assertElementPresent(json,jsonpath)"
$ws.Cells.Item(6,15).Value = "assertLesser(num1,num2)"
$ws.Cells.Item(6,16).Value = "assertPatternNotPresent(pdf,regex)"
$ws.Cells.Item(6,17).Value = "saveResult(db,sql,output)"
$ws.Cells.Item(6,18).Value = "flushDb(profile)"
$ws.Cells.Item(6,21).Value = "sftpDelete(var,profile,remote)"
$ws.Cells.Item(6,23).Value = "assertAttributeNotPresent(locator,attrName)"
$ws.Cells.Item(6,24).Value = "replyCancel(text)"
$ws.Cells.Item(6,25).Value = "deleteAll()"
$ws.Cells.Item(6,26).Value = "head(url,var)"
$ws.Cells.Item(6,27).Value = "patch(url,body,output)"
$ws.Cells.Item(6,28).Value = "assertValue(xml,xpath,expected)"
$ws.Cells.Item(7,1).Value = "excel"
$ws.Cells.Item(7,2).Value = "list(var,profile,remotePath)"
$ws.Cells.Item(7,4).Value = "assertCount(text,regex,expects)"
$ws.Cells.Item(7,6).Value = "assertHierCells(matchBy,column,expected,nestedOnly)"
$ws.Cells.Item(7,7).Value = "json(file,worksheet,range,header,output)"
$ws.Cells.Item(7,10).Value = "copyFiles(source,target)"
$ws.Cells.Item(7,12).Value = "assertEqual(expected,actual)"
$ws.Cells.Item(7,15).Value = "assertLesserOrEqual(num1,num2)"
$ws.Cells.Item(7,16).Value = "assertPatternPresent(pdf,regex)"
$ws.Cells.Item(7,17).Value = "saveResults(db,sqls,outputDir)"
$ws.Cells.Item(7,18).Value = "rename(profile,current,new)"
$ws.Cells.Item(7,21).Value = "sftpList(var,profile,remote)"
$ws.Cells.Item(7,23).Value = "assertAttributePresent(locator,attrName)"
$ws.Cells.Item(7,24).Value = "replyOK(text)"
$ws.Cells.Item(7,25).Value = "save(var,name)"
$ws.Cells.Item(7,26).Value = "header(name,value)"
$ws.Cells.Item(7,27).Value = "post(url,body,output)"
$ws.Cells.Item(7,28).Value = "assertValues(xml,xpath,array,exactOrder)"
$ws.Cells.Item(8,1).Value = "external"
$ws.Cells.Item(8,2).Value = "moveFrom(var,profile,remote,local)"
$ws.Cells.Item(8,4).Value = "assertEmpty(text)"
$ws.Cells.Item(8,6).Value = "assertHierRow(matchBy,expected)"
$ws.Cells.Item(8,7).Value = "saveData(var,file,worksheet,range)"
$ws.Cells.Item(8,10).Value = "count(var,path,pattern)"
$ws.Cells.Item(8,12).Value = "assertValue(json,jsonpath,expected)"
$ws.Cells.Item(8,15).Value = "average(var,array)"
$ws.Cells.Item(8,16).Value = "assertTextArray(pdf,textArray,ordered)"
$ws.Cells.Item(8,18).Value = "set(profile,key,value)"
$ws.Cells.Item(8,21).Value = "sftpMoveFrom(var,profile,remote,local)"
$ws.Cells.Item(8,23).Value = "assertChecked(locator)"
$ws.Cells.Item(8,24).Value = "storeText(var)"
$ws.Cells.Item(8,25).Value = "saveAll(var)"
$ws.Cells.Item(8,26).Value = "headerByVar(name,var)"
$ws.Cells.Item(8,27).Value = "put(url,body,output)"
$ws.Cells.Item(8,28).Value = "assertWellformed(xml)"
$ws.Cells.Item(9,1).Value = "image"
$ws.Cells.Item(9,2).Value = "moveTo(var,profile,local,remote)"
$ws.Cells.Item(9,4).Value = "assertEndsWith(text,suffix)"
$ws.Cells.Item(9,6).Value = "assertListCount(count)"
$ws.Cells.Item(9,7).Value = "saveRange(var,file,worksheet,range)"
$ws.Cells.Item(9,10).Value = "deleteFiles(location,recursive)"
$ws.Cells.Item(9,12).Value = "assertValues(json,jsonpath,array,exactOrder)"
$ws.Cells.Item(9,15).Value = "ceiling(var)"
$ws.Cells.Item(9,16).Value = "assertTextNotPresent(pdf,text)"
$ws.Cells.Item(9,18).Value = "store(var,profile,key)"
$ws.Cells.Item(9,21).Value = "sftpMoveTo(var,profile,local,remote)"
$ws.Cells.Item(9,23).Value = "assertContainCount(locator,text,count)"
$ws.Cells.Item(9,26).Value = "jwtParse(var,token,key)"
$ws.Cells.Item(9,28).Value = "beautify(xml,var)"
$ws.Cells.Item(10,1).Value = "io"
$ws.Cells.Item(10,4).Value = "assertEqual(expected,actual)"
$ws.Cells.Item(10,6).Value = "assertLocatorNotPresent(locator)"
$ws.Cells.Item(10,7).Value = "setPassword(file,password)"
$ws.Cells.Item(10,10).Value = "filter(source,target,matchPattern)"
$ws.Cells.Item(10,12).Value = "assertWellformed(json)"
$ws.Cells.Item(10,15).Value = "decrement(var,amount)"
$ws.Cells.Item(10,16).Value = "assertTextPresent(pdf,text)"
$ws.Cells.Item(10,18).Value = "storeKeys(var,profile,keyPattern)"
$ws.Cells.Item(10,23).Value = "assertCssNotPresent(locator,property)"
$ws.Cells.Item(10,26).Value = "jwtSignHS256(var,payload,key)"
$ws.Cells.Item(10,28).Value = "minify(xml,var)"
$ws.Cells.Item(11,1).Value = "jms"
$ws.Cells.Item(11,4).Value = "assertNotContains(text,substring)"
$ws.Cells.Item(11,6).Value = "assertLocatorPresent(locator)"
$ws.Cells.Item(11,7).Value = "write(file,worksheet,startCell,data)"
$ws.Cells.Item(11,10).Value = "makeDirectory(source)"
$ws.Cells.Item(11,12).Value = "beautify(json,var)"
$ws.Cells.Item(11,15).Value = "floor(var)"
$ws.Cells.Item(11,16).Value = "count(pdf,text,var)"
$ws.Cells.Item(11,23).Value = "assertCssPresent(locator,property,value)"
$ws.Cells.Item(11,26).Value = "oauth(var,url,auth)"
$ws.Cells.Item(11,28).Value = "storeCount(xml,xpath,var)"
$ws.Cells.Item(12,1).Value = "json"
$ws.Cells.Item(12,4).Value = "assertNotEmpty(text)"
$ws.Cells.Item(12,6).Value = "assertMenuEnabled(menu)"
$ws.Cells.Item(12,7).Value = "writeAcross(file,worksheet,startCell,array)"
$ws.Cells.Item(12,10).Value = "moveFiles(source,target)"
$ws.Cells.Item(12,12).Value = "fromCsv(csv,header,jsonFile)"
$ws.Cells.Item(12,15).Value = "increment(var,amount)"
$ws.Cells.Item(12,16).Value = "saveAsPages(pdf,destination)"
$ws.Cells.Item(12,23).Value = "assertElementByAttributes(nameValues)"
$ws.Cells.Item(12,26).Value = "patch(url,body,var)"
$ws.Cells.Item(12,28).Value = "storeValue(xml,xpath,var)"
$ws.Cells.Item(13,1).Value = "macro"
$ws.Cells.Item(13,4).Value = "assertNotEqual(expected,actual)"
$ws.Cells.Item(13,6).Value = "assertModalDialogNotPresent()"
$ws.Cells.Item(13,7).Value = "writeDown(file,worksheet,startCell,array)"
$ws.Cells.Item(13,10).Value = "readFile(var,file)"
$ws.Cells.Item(13,12).Value = "minify(json,var)"
$ws.Cells.Item(13,15).Value = "max(var,array)"
$ws.Cells.Item(13,16).Value = "saveAsText(pdf,destination)"
$ws.Cells.Item(13,23).Value = "assertElementByText(locator,text)"
$ws.Cells.Item(13,26).Value = "post(url,body,var)"
$ws.Cells.Item(13,28).Value = "storeValues(xml,xpath,var)"
$ws.Cells.Item(14,1).Value = "mail"
$ws.Cells.Item(14,4).Value = "assertStartsWith(text,prefix)"
$ws.Cells.Item(14,6).Value = "assertModalDialogPresent()"
$ws.Cells.Item(14,7).Value = "writeVar(var,file,worksheet,startCell)"
$ws.Cells.Item(14,10).Value = "readProperty(var,file,property)"
$ws.Cells.Item(14,12).Value = "storeCount(json,jsonpath,var)"
$ws.Cells.Item(14,15).Value = "min(var,array)"
$ws.Cells.Item(14,16).Value = "saveFormValues(pdf,var,pageAndLineStartEnd,strategy)"
$ws.Cells.Item(14,23).Value = "assertElementCount(locator,count)"
$ws.Cells.Item(14,26).Value = "put(url,body,var)"
$ws.Cells.Item(15,1).Value = "number"
$ws.Cells.Item(15,4).Value = "assertTextOrder(var,descending)"
$ws.Cells.Item(15,6).Value = "assertModalDialogTitle(title)"
$ws.Cells.Item(15,10).Value = "rename(target,newName)"
$ws.Cells.Item(15,12).Value = "storeValue(json,jsonpath,var)"
$ws.Cells.Item(15,15).Value = "round(var,closestDigit)"
$ws.Cells.Item(15,16).Value = "saveMetadata(pdf,var)"
$ws.Cells.Item(15,23).Value = "assertElementNotPresent(locator)"
$ws.Cells.Item(15,26).Value = "saveResponsePayload(var,file,append)"
$ws.Cells.Item(16,1).Value = "pdf"
$ws.Cells.Item(16,4).Value = "assertVarNotPresent(var)"
$ws.Cells.Item(16,6).Value = "assertModalDialogTitleByLocator(locator,title)"
$ws.Cells.Item(16,10).Value = "saveDiff(var,expected,actual)"
$ws.Cells.Item(16,12).Value = "storeValues(json,jsonpath,var)"
$ws.Cells.Item(16,16).Value = "saveToVar(pdf,var)"
$ws.Cells.Item(16,23).Value = "assertElementPresent(locator)"
$ws.Cells.Item(16,26).Value = "soap(action,url,payload,var)"
$ws.Cells.Item(17,1).Value = "rdbms"
$ws.Cells.Item(17,4).Value = "assertVarPresent(var)"
$ws.Cells.Item(17,6).Value = "assertNotChecked(name)"
$ws.Cells.Item(17,10).Value = "saveFileMeta(var,file)"
$ws.Cells.Item(17,23).Value = "assertFocus(locator)"
$ws.Cells.Item(17,26).Value = "upload(url,body,fileParams,var)"
$ws.Cells.Item(18,1).Value = "redis"
$ws.Cells.Item(18,4).Value = "clear(vars)"
$ws.Cells.Item(18,6).Value = "assertSelected(name,text)"
$ws.Cells.Item(18,10).Value = "saveMatches(var,path,filePattern)"
$ws.Cells.Item(18,23).Value = "assertFrameCount(count)"
$ws.Cells.Item(19,1).Value = "sms"
$ws.Cells.Item(19,4).Value = "failImmediate(text)"
$ws.Cells.Item(19,6).Value = "assertTableCell(row,column,contains)"
$ws.Cells.Item(19,10).Value = "searchAndReplace(file,config,saveAs)"
$ws.Cells.Item(19,23).Value = "assertFramePresent(frameName)"
$ws.Cells.Item(20,1).Value = "sound"
$ws.Cells.Item(20,4).Value = "incrementChar(var,amount,config)"
$ws.Cells.Item(20,6).Value = "assertTableColumnContains(column,contains)"
$ws.Cells.Item(20,10).Value = "unzip(zipFile,target)"
$ws.Cells.Item(20,23).Value = "assertIECompatMode()"
$ws.Cells.Item(21,1).Value = "ssh"
$ws.Cells.Item(21,4).Value = "macro(file,sheet,name)"
$ws.Cells.Item(21,6).Value = "assertTableContains(contains)"
$ws.Cells.Item(21,10).Value = "validate(var,profile,inputFile)"
$ws.Cells.Item(21,23).Value = "assertIENavtiveMode()"
$ws.Cells.Item(22,1).Value = "step"
$ws.Cells.Item(22,4).Value = "prependText(var,prependWith)"
$ws.Cells.Item(22,6).Value = "assertTableRowContains(row,contains)"
$ws.Cells.Item(22,10).Value = "writeFile(file,content,append)"
$ws.Cells.Item(22,23).Value = "assertLinkByLabel(label)"
$ws.Cells.Item(23,1).Value = "web"
$ws.Cells.Item(23,4).Value = "repeatUntil(steps,maxWaitMs)"
$ws.Cells.Item(23,6).Value = "assertText(name,expected)"
$ws.Cells.Item(23,10).Value = "writeFileAsIs(file,content,append)"
$ws.Cells.Item(23,23).Value = "assertNotChecked(locator)"
$ws.Cells.Item(24,1).Value = "webalert"
$ws.Cells.Item(24,4).Value = "save(var,value)"
$ws.Cells.Item(24,6).Value = "assertWindowTitleContains(contains)"
$ws.Cells.Item(24,10).Value = "writeProperty(file,property,value)"
$ws.Cells.Item(24,23).Value = "assertNotFocus(locator)"
$ws.Cells.Item(25,1).Value = "webcookie"
$ws.Cells.Item(25,4).Value = "saveCount(text,regex,saveVar)"
$ws.Cells.Item(25,6).Value = "clear(locator)"
$ws.Cells.Item(25,10).Value = "zip(filePattern,zipFile)"
$ws.Cells.Item(25,23).Value = "assertNotText(locator,text)"
$ws.Cells.Item(26,1).Value = "ws"
$ws.Cells.Item(26,4).Value = "saveMatches(text,regex,saveVar)"
$ws.Cells.Item(26,6).Value = "clearCombo(name)"
$ws.Cells.Item(26,23).Value = "assertNotVisible(locator)"
$ws.Cells.Item(27,1).Value = "ws.async"
$ws.Cells.Item(27,4).Value = "saveReplace(text,regex,replace,saveVar)"
$ws.Cells.Item(27,6).Value = "clearModalDialog(var,button)"
$ws.Cells.Item(27,23).Value = "assertOneMatch(locator)"
$ws.Cells.Item(28,1).Value = "xml"
$ws.Cells.Item(28,4).Value = "saveVariablesByPrefix(var,prefix)"
$ws.Cells.Item(28,6).Value = "clearTextArea(name)"
$ws.Cells.Item(28,23).Value = "assertScrollbarHNotPresent(locator)"
$ws.Cells.Item(29,4).Value = "saveVariablesByRegex(var,regex)"
$ws.Cells.Item(29,6).Value = "clearTextBox(name)"
$ws.Cells.Item(29,23).Value = "assertScrollbarHPresent(locator)"
$ws.Cells.Item(30,4).Value = "section(steps)"
$ws.Cells.Item(30,6).Value = "clickButton(name)"
$ws.Cells.Item(30,23).Value = "assertScrollbarVNotPresent(locator)"
$ws.Cells.Item(31,4).Value = "split(text,delim,saveVar)"
$ws.Cells.Item(31,6).Value = "clickByLocator(locator)"
$ws.Cells.Item(31,23).Value = "assertScrollbarVPresent(locator)"
$ws.Cells.Item(32,4).Value = "startRecording()"
$ws.Cells.Item(32,6).Value = "clickCheckBox(name)"
$ws.Cells.Item(32,23).Value = "assertTable(locator,row,column,text)"
$ws.Cells.Item(33,4).Value = "stopRecording()"
$ws.Cells.Item(33,6).Value = "clickExplorerBar(group,item)"
$ws.Cells.Item(33,23).Value = "assertText(locator,text)"
$ws.Cells.Item(34,4).Value = "substringAfter(text,delim,saveVar)"
$ws.Cells.Item(34,6).Value = "clickFirstMatchRow(nameValues)"
$ws.Cells.Item(34,23).Value = "assertTextContains(locator,text)"
$ws.Cells.Item(35,4).Value = "substringBefore(text,delim,saveVar)"
$ws.Cells.Item(35,6).Value = "clickFirstMatchedList(contains)"
$ws.Cells.Item(35,23).Value = "assertTextCount(locator,text,count)"
$ws.Cells.Item(36,4).Value = "substringBetween(text,start,end,saveVar)"
$ws.Cells.Item(36,6).Value = "clickIcon(label)"
$ws.Cells.Item(36,23).Value = "assertTextList(locator,list,ignoreOrder)"
$ws.Cells.Item(37,4).Value = "verbose(text)"
$ws.Cells.Item(37,6).Value = "clickList(row)"
$ws.Cells.Item(37,23).Value = "assertTextMatches(text,minMatch,scrollTo)"
$ws.Cells.Item(38,4).Value = "waitFor(waitMs)"
$ws.Cells.Item(38,6).Value = "clickMenu(menu)"
$ws.Cells.Item(38,23).Value = "assertTextNotPresent(text)"
$ws.Cells.Item(39,6).Value = "clickOffset(locator,xOffset,yOffset)"
$ws.Cells.Item(39,23).Value = "assertTextOrder(locator,descending)"
$ws.Cells.Item(40,6).Value = "clickRadio(name)"
$ws.Cells.Item(40,23).Value = "assertTextPresent(text)"
$ws.Cells.Item(41,6).Value = "clickTab(group,name)"
$ws.Cells.Item(41,23).Value = "assertTitle(text)"
$ws.Cells.Item(42,6).Value = "clickTableCell(row,column)"
$ws.Cells.Item(42,23).Value = "assertValue(locator,value)"
$ws.Cells.Item(43,6).Value = "clickTableRow(row)"
$ws.Cells.Item(43,23).Value = "assertValueOrder(locator,descending)"
$ws.Cells.Item(44,6).Value = "clickTextPane(name,criteria)"
$ws.Cells.Item(44,23).Value = "assertVisible(locator)"
$ws.Cells.Item(45,6).Value = "clickTextPaneRow(var,index)"
$ws.Cells.Item(45,23).Value = "checkAll(locator)"
$ws.Cells.Item(46,6).Value = "closeApplication()"
$ws.Cells.Item(46,23).Value = "clearLocalStorage()"
$ws.Cells.Item(47,6).Value = "collapseHierTable()"
$ws.Cells.Item(47,23).Value = "click(locator)"
$ws.Cells.Item(48,6).Value = "editCurrentRow(nameValues)"
$ws.Cells.Item(48,23).Value = "clickAndWait(locator,waitMs)"
$ws.Cells.Item(49,6).Value = "editHierCells(var,matchBy,nameValues)"
$ws.Cells.Item(49,23).Value = "clickByLabel(label)"
$ws.Cells.Item(50,6).Value = "editTableCells(row,nameValues)"
$ws.Cells.Item(50,23).Value = "clickByLabelAndWait(label,waitMs)"
$ws.Cells.Item(51,6).Value = "getRowCount(var)"
$ws.Cells.Item(51,23).Value = "clickOffset(locator,x,y)"
$ws.Cells.Item(52,6).Value = "hideExplorerBar()"
$ws.Cells.Item(52,23).Value = "clickWithKeys(locator,keys)"
$ws.Cells.Item(53,6).Value = "login(form,username,password)"
$ws.Cells.Item(53,23).Value = "close()"
$ws.Cells.Item(54,6).Value = "maximize()"
$ws.Cells.Item(54,23).Value = "closeAll()"
$ws.Cells.Item(55,6).Value = "minimize()"
$ws.Cells.Item(55,23).Value = "deselect(locator,text)"
$ws.Cells.Item(56,6).Value = "resize(width,height)"
$ws.Cells.Item(56,23).Value = "deselectMulti(locator,array)"
$ws.Cells.Item(57,6).Value = "saveAllTableRows(var)"
$ws.Cells.Item(57,23).Value = "dismissInvalidCert()"
$ws.Cells.Item(58,6).Value = "saveAttributeByLocator(var,locator,attribute)"
$ws.Cells.Item(58,23).Value = "dismissInvalidCertPopup()"
$ws.Cells.Item(59,6).Value = "saveElementCount(var,name)"
$ws.Cells.Item(59,23).Value = "doubleClick(locator)"
$ws.Cells.Item(60,6).Value = "saveFirstListData(var,contains)"
$ws.Cells.Item(60,23).Value = "doubleClickAndWait(locator,waitMs)"
$ws.Cells.Item(61,6).Value = "saveFirstMatchedListIndex(var,contains)"
$ws.Cells.Item(61,23).Value = "doubleClickByLabel(label)"
$ws.Cells.Item(62,6).Value = "saveHierCells(var,matchBy,column,nestedOnly)"
$ws.Cells.Item(62,23).Value = "doubleClickByLabelAndWait(label,waitMs)"
$ws.Cells.Item(63,6).Value = "saveHierRow(var,matchBy)"
$ws.Cells.Item(63,23).Value = "dragAndDrop(fromLocator,toLocator)"
$ws.Cells.Item(64,6).Value = "saveListData(var,contains)"
$ws.Cells.Item(64,23).Value = "dragTo(fromLocator,xOffset,yOffset)"
$ws.Cells.Item(65,6).Value = "saveLocatorCount(var,locator)"
$ws.Cells.Item(65,23).Value = "editLocalStorage(key,value)"
$ws.Cells.Item(66,6).Value = "saveModalDialogText(var)"
$ws.Cells.Item(66,23).Value = "executeScript(var,script)"
$ws.Cells.Item(67,6).Value = "saveModalDialogTextByLocator(var,locater)"
$ws.Cells.Item(67,23).Value = "focus(locator)"
$ws.Cells.Item(68,6).Value = "saveProcessId(var,locator)"
$ws.Cells.Item(68,23).Value = "goBack()"
$ws.Cells.Item(69,6).Value = "saveRowCount(var)"
$ws.Cells.Item(69,23).Value = "goBackAndWait()"
$ws.Cells.Item(70,6).Value = "saveTableRows(var,contains)"
$ws.Cells.Item(70,23).Value = "maximizeWindow()"
$ws.Cells.Item(71,6).Value = "saveTableRowsRange(var,beginRow,endRow)"
$ws.Cells.Item(71,23).Value = "mouseOver(locator)"
$ws.Cells.Item(72,6).Value = "saveText(var,name)"
$ws.Cells.Item(72,23).Value = "open(url)"
$ws.Cells.Item(73,6).Value = "saveTextPane(var,name,criteria)"
$ws.Cells.Item(73,23).Value = "openAndWait(url,waitMs)"
$ws.Cells.Item(74,6).Value = "saveWindowTitle(var)"
$ws.Cells.Item(74,23).Value = "openHttpBasic(url,username,password)"
$ws.Cells.Item(75,6).Value = "scanTable(var,name)"
$ws.Cells.Item(75,23).Value = "openIgnoreTimeout(url)"
$ws.Cells.Item(76,6).Value = "selectCombo(name,text)"
$ws.Cells.Item(76,23).Value = "refresh()"
$ws.Cells.Item(77,6).Value = "sendKeysToTextBox(name,text1,text2,text3,text4)"
$ws.Cells.Item(77,23).Value = "refreshAndWait()"
$ws.Cells.Item(78,6).Value = "showExplorerBar()"
$ws.Cells.Item(78,23).Value = "resizeWindow(width,height)"
$ws.Cells.Item(79,6).Value = "toggleExplorerBar()"
$ws.Cells.Item(79,23).Value = "saveAllWindowIds(var)"
$ws.Cells.Item(80,6).Value = "typeAppendTextArea(name,text1,text2,text3,text4)"
$ws.Cells.Item(80,23).Value = "saveAllWindowNames(var)"
$ws.Cells.Item(81,6).Value = "typeAppendTextBox(name,text1,text2,text3,text4)"
$ws.Cells.Item(81,23).Value = "saveAttribute(var,locator,attrName)"
$ws.Cells.Item(82,6).Value = "typeByLocator(locator,text)"
$ws.Cells.Item(82,23).Value = "saveAttributeList(var,locator,attrName)"
$ws.Cells.Item(83,6).Value = "typeTextArea(name,text1,text2,text3,text4)"
$ws.Cells.Item(83,23).Value = "saveCount(var,locator)"
$ws.Cells.Item(84,6).Value = "typeTextBox(name,text1,text2,text3,text4)"
$ws.Cells.Item(84,23).Value = "saveDivsAsCsv(headers,rows,cells,nextPage,file)"
$ws.Cells.Item(85,6).Value = "useApp(appId)"
$ws.Cells.Item(85,23).Value = "saveElement(var,locator)"
$ws.Cells.Item(86,6).Value = "useForm(formName)"
$ws.Cells.Item(86,23).Value = "saveElements(var,locator)"
$ws.Cells.Item(87,6).Value = "useHierTable(var,name)"
$ws.Cells.Item(87,23).Value = "saveLocalStorage(var,key)"
$ws.Cells.Item(88,6).Value = "useList(var,name)"
$ws.Cells.Item(88,23).Value = "saveLocation(var)"
$ws.Cells.Item(89,6).Value = "useTable(var,name)"
$ws.Cells.Item(89,23).Value = "savePageAs(var,sessionIdName,url)"
$ws.Cells.Item(90,6).Value = "useTableRow(var,row)"
$ws.Cells.Item(90,23).Value = "savePageAsFile(sessionIdName,url,file)"
$ws.Cells.Item(91,6).Value = "waitFor(name,maxWaitMs)"
$ws.Cells.Item(91,23).Value = "saveTableAsCsv(locator,nextPageLocator,file)"
$ws.Cells.Item(92,6).Value = "waitForLocator(locator,maxWaitMs)"
$ws.Cells.Item(92,23).Value = "saveText(var,locator)"
$ws.Cells.Item(93,23).Value = "saveTextArray(var,locator)"
$ws.Cells.Item(94,23).Value = "saveTextSubstringAfter(var,locator,delim)"
$ws.Cells.Item(95,23).Value = "saveTextSubstringBefore(var,locator,delim)"
$ws.Cells.Item(96,23).Value = "saveTextSubstringBetween(var,locator,start,end)"
$ws.Cells.Item(97,23).Value = "saveValue(var,locator)"
$ws.Cells.Item(98,23).Value = "scrollLeft(locator,pixel)"
$ws.Cells.Item(99,23).Value = "scrollRight(locator,pixel)"
$ws.Cells.Item(100,23).Value = "scrollTo(locator)"
$ws.Cells.Item(101,23).Value = "select(locator,text)"
$ws.Cells.Item(102,23).Value = "selectFrame(locator)"
$ws.Cells.Item(103,23).Value = "selectMulti(locator,array)"
$ws.Cells.Item(104,23).Value = "selectMultiOptions(locator)"
$ws.Cells.Item(105,23).Value = "selectText(locator)"
$ws.Cells.Item(106,23).Value = "selectWindow(winId)"
$ws.Cells.Item(107,23).Value = "selectWindowAndWait(winId,waitMs)"
$ws.Cells.Item(108,23).Value = "selectWindowByIndex(index)"
$ws.Cells.Item(109,23).Value = "selectWindowByIndexAndWait(index,waitMs)"
$ws.Cells.Item(110,23).Value = "toggleSelections(locator)"
$ws.Cells.Item(111,23).Value = "type(locator,value)"
$ws.Cells.Item(112,23).Value = "typeKeys(locator,value)"
$ws.Cells.Item(113,23).Value = "uncheckAll(locator)"
$ws.Cells.Item(114,23).Value = "unselectAllText()"
$ws.Cells.Item(115,23).Value = "upload(fieldLocator,file)"
$ws.Cells.Item(116,23).Value = "verifyContainText(locator,text)"
$ws.Cells.Item(117,23).Value = "verifyText(locator,text)"
$ws.Cells.Item(118,23).Value = "wait(waitMs)"
$ws.Cells.Item(119,23).Value = "waitForElementPresent(locator)"
$ws.Cells.Item(120,23).Value = "waitForPopUp(winId,waitMs)"
$ws.Cells.Item(121,23).Value = "waitForTextPresent(text)"
$ws.Cells.Item(122,23).Value = "waitForTitle(text)"
# Update defined names to reflect the shifted / extended ranges
$wb.Names.Item("external").RefersTo = "='#system'!`$H`$2:`$H`$4"
$wb.Names.Item("mail").RefersTo = "='#system'!`$N`$2:`$N`$2"
$wb.Names.Item("number").RefersTo = "='#system'!`$O`$2:`$O`$15"
$wb.Names.Item("pdf").RefersTo = "='#system'!`$P`$2:`$P`$16"
$wb.Names.Item("rdbms").RefersTo = "='#system'!`$Q`$2:`$Q`$7"
$wb.Names.Item("redis").RefersTo = "='#system'!`$R`$2:`$R`$10"
$wb.Names.Item("sms").RefersTo = "='#system'!`$S`$2:`$S`$2"
$wb.Names.Item("sound").RefersTo = "='#system'!`$T`$2:`$T`$5"
$wb.Names.Item("ssh").RefersTo = "='#system'!`$U`$2:`$U`$9"
$wb.Names.Item("step").RefersTo = "='#system'!`$V`$2:`$V`$4"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$28"
$wb.Names.Item("web").RefersTo = "='#system'!`$W`$2:`$W`$122"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$X`$2:`$X`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$Y`$2:`$Y`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$Z`$2:`$Z`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AB`$2:`$AB`$13"
$wb.Names.Add("macro", "='#system'!`$M`$2:`$M`$4")
